$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '28.253.25'
$ws.Range('E2').Value = '  -0.56%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.804.49'
$ws.Range('E3').Value = '  -0.75%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.13%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '314.12'
$ws.Range('E5').Value = '  -0.41%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.003'
$ws.Range('E6').Value = '  +0.15%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5280'
$ws.Range('E7').Value = '  +3.30%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3826'
$ws.Range('E8').Value = '  -3.26%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.08027'
$ws.Range('E9').Value = '  -0.69%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '41.40'
$ws.Range('E10').Value = '  -0.63%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.101'
$ws.Range('E11').Value = '  -0.54%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '6.327'
$ws.Range('E12').Value = '  +1.03%  '

# Row 13
$ws.Range('E13').Value = '  +0.14%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '20.62'
$ws.Range('E14').Value = '  -1.77%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '1.808.06'
$ws.Range('E15').Value = '  -0.44%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '7.322'
$ws.Range('E16').Value = '  -2.41%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '92.21'
$ws.Range('E17').Value = '  -0.46%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.00001097'
$ws.Range('E18').Value = '  -3.77%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06607'
$ws.Range('E19').Value = '  -0.46%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.003'
$ws.Range('E20').Value = '  +0.14%  '

# Row 21
$ws.Range('E21').Value = '  -1.76%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.968'
$ws.Range('E22').Value = '  -1.99%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '28.306.89'
$ws.Range('E23').Value = '  -0.47%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.16'
$ws.Range('E24').Value = '  -1.00%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.260'
$ws.Range('E25').Value = '  -0.32%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '160.72'
$ws.Range('E26').Value = '  +3.72%  '

# Row 27
$ws.Range('E27').Value = '  -3.13%  '

# Row 28
$ws.Range('E28').Value = '  -0.88%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.358'
$ws.Range('E29').Value = '  -1.83%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '123.39'
$ws.Range('E30').Value = '  -1.95%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.1085'
$ws.Range('E31').Value = '  -1.44%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.057'
$ws.Range('E32').Value = '  -4.11%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.688'
$ws.Range('E33').Value = '  +0.91%  '

# Row 34
$ws.Range('E34').Value = '  -3.59%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.07248'
$ws.Range('E35').Value = '  +3.10%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '12.38'
$ws.Range('E36').Value = '  +9.66%  '

# Row 37
$ws.Range('E37').Value = '  -0.28%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2154'
$ws.Range('E38').Value = '  -3.22%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.112'
$ws.Range('E39').Value = '  -2.05%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '8.663'
$ws.Range('E40').Value = '  -1.73%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.6204'
$ws.Range('E41').Value = '  -0.96%  '

# Row 42
$ws.Range('E42').Value = '  -0.67%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.372'
$ws.Range('E43').Value = '  -2.12%  '

# Row 44
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.31'
$ws.Range('E44').Value = '  -1.52%  '

# Row 45
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.6024'
$ws.Range('E45').Value = '  +1.89%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.766'
$ws.Range('E46').Value = '  +0.67%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '127.13'
$ws.Range('E47').Value = '  +1.79%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.216'
$ws.Range('E48').Value = '  +2.42%  '

# Row 49
$ws.Range('E49').Value = '  -2.25%  '

# Row 50
$ws.Range('E50').Value = '  -0.93%  '

# Row 51
$ws.Range('E51').Value = '  -1.62%  '
